$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new end-user row (NAME, EMAIL, EMPLOYEE_CODE)
$ws.Range("A4").Value = "Yakshitha"
$ws.Range("B4").Value = "ykmangalore100@gmail.com"
$ws.Range("C4").Value = "YK001"

# Turn the new email into a mailto hyperlink, same as the existing B2 entry
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:ykmangalore100@gmail.com", "", "", "ykmangalore100@gmail.com")

# Match the formatting already used for the other hyperlinked email (B2):
# Arial 10, non-underlined, blue font instead of Excel's default Hyperlink style
$ws.Range("B4").Font.Underline = $false
$ws.Range("B4").Font.Color = $ws.Range("B2").Font.Color
$ws.Range("B4").Font.Name = "Arial"
$ws.Range("B4").Font.Size = 10

# Move/keep the active selection on C4, matching the saved view state
$ws.Range("C4").Select()
